# New weekly price observation for "Ajo" (Chino, Primera) at Terminal
# Hortofrutícola Agro Chillán. The new record is inserted as the first data
# row (row 93) of this variety/quality block, pushing every subsequent row
# down by one (Excel's normal row-insert behaviour), so the whole block
# effectively keeps its most-recent-first ordering with the new week on top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 93 - everything below (old rows
# 93..186) shifts down to 94..187 automatically, carrying its data/styles
# with it.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new observation. Most fields
# mirror the block's constant columns (mercado, región, categoría, etc.);
# only the date (D) is genuinely new, the rest of this particular record
# duplicates the values of what is now row 94 (price bucket 18000-19000).
$ws.Range("A93").Value = 7
$ws.Range("B93").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C93").Value = "Ñuble"
$ws.Range("D93").Value = 44587
$ws.Range("E93").Value = 16
$ws.Range("F93").Value = 100112003
$ws.Range("G93").Value = "Ajo"
$ws.Range("H93").Value = "Chino"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 60
$ws.Range("K93").Value = 18000
$ws.Range("L93").Value = 19000
$ws.Range("M93").Value = 18500
$ws.Range("N93").Value = "`$/caja 10 kilos"
$ws.Range("O93").Value = "China"
$ws.Range("P93").Value = 1850
$ws.Range("Q93").Value = 10
$ws.Range("R93").Value = "Hortaliza"
